# Fruta / hortaliza, semanal
#
# A new weekly price-record row is inserted into the "Hortaliza, Macroferia
# Regional de Talca - Betarraga" table at row 195 (pushing the former rows
# 195-206 down to 196-207, and growing the used range from A1:R206 to
# A1:R207). The new row carries the same market/category/region metadata as
# its neighbours, with its own date and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 195-206 down to 196-207, leaving a blank row 195 to populate.
$ws.Rows.Item(195).Insert()

# Fill in the newly inserted row 195.
$ws.Cells.Item(195, 1).Value  = 5
$ws.Cells.Item(195, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(195, 3).Value  = "Maule"
$ws.Cells.Item(195, 4).Value  = 44516
$ws.Cells.Item(195, 5).Value  = 7
$ws.Cells.Item(195, 6).Value  = 100114014
$ws.Cells.Item(195, 7).Value  = "Betarraga"
$ws.Cells.Item(195, 8).Value  = "Sin especificar"
$ws.Cells.Item(195, 9).Value  = "Primera"
$ws.Cells.Item(195, 10).Value = 5000
$ws.Cells.Item(195, 11).Value = 550
$ws.Cells.Item(195, 12).Value = 550
$ws.Cells.Item(195, 13).Value = 550
$ws.Cells.Item(195, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(195, 15).Value = "Región del Maule"
$ws.Cells.Item(195, 16).Value = 110
$ws.Cells.Item(195, 17).Value = 5
$ws.Cells.Item(195, 18).Value = "Hortaliza"
